$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 7.781999999999999
$ws.Range("B7").Value = 5.954
$ws.Range("D7").Value = -7.567000000000002
$ws.Range("D15").Value = -8.411
$ws.Range("B16").Value = 4.843
$ws.Range("E16").Value = 16.409
$ws.Range("E19").Value = 16.593
$ws.Range("D21").Value = -8.1
$ws.Range("D22").Value = -7.948
$ws.Range("D23").Value = -7.869999999999999
$ws.Range("B28").Value = 6.034000000000001
$ws.Range("B29").Value = 5.207
$ws.Range("B32").Value = 6.572
$ws.Range("D34").Value = -7.917999999999999
$ws.Range("E36").Value = 16.691
$ws.Range("B40").Value = 9.223000000000001
$ws.Range("D43").Value = -7.712000000000001
$ws.Range("D45").Value = -7.532999999999999
$ws.Range("E46").Value = 16.942
$ws.Range("D50").Value = -8.135
$ws.Range("E50").Value = 16.605
$ws.Range("D51").Value = -8.384
$ws.Range("B52").Value = 4.944000000000001
$ws.Range("B57").Value = 5.091
$ws.Range("B66").Value = 5.773
$ws.Range("D66").Value = -7.561000000000002
$ws.Range("D67").Value = -6.783999999999999
$ws.Range("D79").Value = -7.614
$ws.Range("D84").Value = -8.132000000000001
$ws.Range("D92").Value = -6.544
$ws.Range("E95").Value = 17.483
$ws.Range("D97").Value = -8.102
$ws.Range("E97").Value = 16.845
$ws.Range("B100").Value = 5.931
